# Auto-generated Excel COM-interop script to apply kraw_summaries.xlsx updates
# per commit: "Update buddingtonite analysis using XPP matrix correction method"
$wb = $excel.ActiveWorkbook

# --- Sheet: Na ---
$ws = $wb.Worksheets.Item("Na")
$ws.Range("H2").Value = 0.72
$ws.Range("J2").Value = 0.72
$ws.Range("K2").Value = 27.62
$ws.Range("L2").Value = 0.72
$ws.Range("M2").Value = 27.62
$ws.Range("K3").Value = 13.25
$ws.Range("M3").Value = 13.25
$ws.Range("H4").Value = 1.02
$ws.Range("J4").Value = 1.02
$ws.Range("K4").Value = 16.87
$ws.Range("L4").Value = 1.02
$ws.Range("M4").Value = 16.87
$ws.Range("H5").Value = 1.19
$ws.Range("J5").Value = 1.19
$ws.Range("K5").Value = 17.55
$ws.Range("L5").Value = 1.19
$ws.Range("M5").Value = 17.55
$ws.Range("H6").Value = 0.82
$ws.Range("J6").Value = 0.82
$ws.Range("K6").Value = 23.22
$ws.Range("L6").Value = 0.82
$ws.Range("M6").Value = 23.22

# --- Sheet: Si ---
$ws = $wb.Worksheets.Item("Si")
$ws.Range("H2").Value = 124.05
$ws.Range("J2").Value = 124.05
$ws.Range("K2").Value = 0.27
$ws.Range("L2").Value = 124.05
$ws.Range("M2").Value = 0.27
$ws.Range("H3").Value = 122.78
$ws.Range("J3").Value = 122.78
$ws.Range("L3").Value = 122.78
$ws.Range("H4").Value = 121.77
$ws.Range("J4").Value = 121.77
$ws.Range("L4").Value = 121.77
$ws.Range("H6").Value = 123.11
$ws.Range("J6").Value = 123.11
$ws.Range("L6").Value = 123.11

# --- Sheet: Al ---
$ws = $wb.Worksheets.Item("Al")
$ws.Range("H2").Value = 99.83
$ws.Range("J2").Value = 99.83
$ws.Range("K2").Value = 0.43
$ws.Range("L2").Value = 99.83
$ws.Range("M2").Value = 0.43
$ws.Range("H3").Value = 100.86
$ws.Range("J3").Value = 100.86
$ws.Range("L3").Value = 100.86
$ws.Range("H4").Value = 100.8
$ws.Range("J4").Value = 100.8
$ws.Range("K4").Value = 0.44
$ws.Range("L4").Value = 100.8
$ws.Range("M4").Value = 0.44
$ws.Range("H5").Value = 101.01
$ws.Range("J5").Value = 101.01
$ws.Range("K5").Value = 0.43
$ws.Range("L5").Value = 101.01
$ws.Range("M5").Value = 0.43
$ws.Range("H6").Value = 99.8
$ws.Range("J6").Value = 99.8
$ws.Range("K6").Value = 0.41
$ws.Range("L6").Value = 99.8
$ws.Range("M6").Value = 0.41

# --- Sheet: Mg ---
$ws = $wb.Worksheets.Item("Mg")
$ws.Range("K2").Value = 85.2
$ws.Range("M2").Value = 85.2
$ws.Range("K3").Value = 161.86
$ws.Range("M3").Value = 161.86
$ws.Range("K4").Value = 105.94
$ws.Range("M4").Value = 105.94
$ws.Range("H5").Value = -0.03
$ws.Range("J5").Value = -0.03
$ws.Range("K5").Value = 140.77
$ws.Range("L5").Value = -0.03
$ws.Range("M5").Value = 140.77
$ws.Range("H6").Value = 0.01
$ws.Range("J6").Value = 0.01
$ws.Range("K6").Value = 269.69
$ws.Range("L6").Value = 0.01
$ws.Range("M6").Value = 269.69

# --- Sheet: K ---
$ws = $wb.Worksheets.Item("K")
$ws.Range("H2").Value = 0.23
$ws.Range("J2").Value = 0.23
$ws.Range("K2").Value = 32.29
$ws.Range("L2").Value = 0.23
$ws.Range("M2").Value = 32.29
$ws.Range("H3").Value = 0.36
$ws.Range("J3").Value = 0.36
$ws.Range("K3").Value = 24.34
$ws.Range("L3").Value = 0.36
$ws.Range("M3").Value = 24.34
$ws.Range("H4").Value = 0.1
$ws.Range("J4").Value = 0.1
$ws.Range("K4").Value = 63.82
$ws.Range("L4").Value = 0.1
$ws.Range("M4").Value = 63.82
$ws.Range("H5").Value = 0.27
$ws.Range("J5").Value = 0.27
$ws.Range("K5").Value = 26.72
$ws.Range("L5").Value = 0.27
$ws.Range("M5").Value = 26.72
$ws.Range("H6").Value = 0.25
$ws.Range("J6").Value = 0.25
$ws.Range("K6").Value = 30.27
$ws.Range("L6").Value = 0.25
$ws.Range("M6").Value = 30.27

# --- Sheet: Ca ---
$ws = $wb.Worksheets.Item("Ca")
$ws.Range("K2").Value = 29.27
$ws.Range("M2").Value = 29.27
$ws.Range("K3").Value = 77.34
$ws.Range("M3").Value = 77.34
$ws.Range("K4").Value = 46.75
$ws.Range("M4").Value = 46.75
$ws.Range("H5").Value = 0.02
$ws.Range("J5").Value = 0.02
$ws.Range("K5").Value = 252.55
$ws.Range("L5").Value = 0.02
$ws.Range("M5").Value = 252.55
$ws.Range("H6").Value = 0.04
$ws.Range("J6").Value = 0.04
$ws.Range("K6").Value = 117.93
$ws.Range("L6").Value = 0.04
$ws.Range("M6").Value = 117.93

# --- Sheet: Rb ---
$ws = $wb.Worksheets.Item("Rb")
$ws.Range("H2").Value = 1.76
$ws.Range("J2").Value = 1.76
$ws.Range("K2").Value = 117.97
$ws.Range("L2").Value = 1.76
$ws.Range("M2").Value = 117.97
$ws.Range("H3").Value = -1.25
$ws.Range("J3").Value = -1.25
$ws.Range("K3").Value = 183.91
$ws.Range("L3").Value = -1.25
$ws.Range("M3").Value = 183.91
$ws.Range("H4").Value = 2.68
$ws.Range("J4").Value = 2.68
$ws.Range("K4").Value = 74.76000000000001
$ws.Range("L4").Value = 2.68
$ws.Range("M4").Value = 74.76000000000001
$ws.Range("H5").Value = 3.58
$ws.Range("J5").Value = 3.58
$ws.Range("K5").Value = 55.71
$ws.Range("L5").Value = 3.58
$ws.Range("M5").Value = 55.71
$ws.Range("H6").Value = 1.84
$ws.Range("J6").Value = 1.84
$ws.Range("K6").Value = 106.1
$ws.Range("L6").Value = 1.84
$ws.Range("M6").Value = 106.1

# --- Sheet: Mo ---
$ws = $wb.Worksheets.Item("Mo")
$ws.Range("K2").Value = 113.96
$ws.Range("M2").Value = 113.96
$ws.Range("K3").Value = 61.9
$ws.Range("M3").Value = 61.9
$ws.Range("K4").Value = 64.2
$ws.Range("M4").Value = 64.2
$ws.Range("K5").Value = 139.37
$ws.Range("M5").Value = 139.37
$ws.Range("K6").Value = 137.15
$ws.Range("M6").Value = 137.15

# --- Sheet: N ---
$ws = $wb.Worksheets.Item("N")
$ws.Range("H2").Value = 32.5
$ws.Range("J2").Value = 26.62
$ws.Range("K2").Value = 3.61
$ws.Range("L2").Value = 32.5
$ws.Range("M2").Value = 3.94
$ws.Range("H3").Value = 33.85
$ws.Range("J3").Value = 27.72
$ws.Range("K3").Value = 3.18
$ws.Range("L3").Value = 33.85
$ws.Range("M3").Value = 3.55
$ws.Range("H4").Value = 31.37
$ws.Range("J4").Value = 25.69
$ws.Range("K4").Value = 3.42
$ws.Range("L4").Value = 31.37
$ws.Range("M4").Value = 3.77
$ws.Range("H5").Value = 30.47
$ws.Range("J5").Value = 24.95
$ws.Range("K5").Value = 3.68
$ws.Range("L5").Value = 30.47
$ws.Range("M5").Value = 4.01
$ws.Range("H6").Value = 30.71
$ws.Range("J6").Value = 25.15
$ws.Range("K6").Value = 3.45
$ws.Range("L6").Value = 30.71
$ws.Range("M6").Value = 3.8

$wb.Save()
